$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 71-74 ---

# Row 71: date changes from 44194 to 44568 (2022-01-07)
$ws.Range("D71").Value = 44568

# Row 72: date changes from 44194 to 44568
$ws.Range("D72").Value = 44568

# Row 73: date changes from 44194 to 44568; volume J changes from 400 to 500
$ws.Range("D73").Value = 44568
$ws.Range("J73").Value = 500

# Row 74: date changes from 44272 to 44194; quality changes Primera -> Extra;
# volume/prices updated accordingly
$ws.Range("D74").Value = 44194
$ws.Range("I74").Value = "Extra"
$ws.Range("J74").Value = 400
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = 3000
$ws.Range("P74").Value = 3000

# --- Append new rows 75-78 ---

# Row 75
$ws.Range("A75").Value = 11
$ws.Range("B75").Value = "Vega Monumental Concepción"
$ws.Range("C75").Value = "Bíobío"
$ws.Range("D75").Value = 44194
$ws.Range("E75").Value = 8
$ws.Range("F75").Value = 100112028
$ws.Range("G75").Value = "Sandia"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 500
$ws.Range("K75").Value = 2500
$ws.Range("L75").Value = 2500
$ws.Range("M75").Value = 2500
$ws.Range("N75").Value = "$/unidad"
$ws.Range("O75").Value = "Región de O'Higgins"
$ws.Range("P75").Value = 2500
$ws.Range("Q75").Value = 1
$ws.Range("R75").Value = "Hortaliza"

# Row 76
$ws.Range("A76").Value = 11
$ws.Range("B76").Value = "Vega Monumental Concepción"
$ws.Range("C76").Value = "Bíobío"
$ws.Range("D76").Value = 44194
$ws.Range("E76").Value = 8
$ws.Range("F76").Value = 100112028
$ws.Range("G76").Value = "Sandia"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Segunda"
$ws.Range("J76").Value = 400
$ws.Range("K76").Value = 2000
$ws.Range("L76").Value = 2000
$ws.Range("M76").Value = 2000
$ws.Range("N76").Value = "$/unidad"
$ws.Range("O76").Value = "Región de O'Higgins"
$ws.Range("P76").Value = 2000
$ws.Range("Q76").Value = 1
$ws.Range("R76").Value = "Hortaliza"

# Row 77
$ws.Range("A77").Value = 11
$ws.Range("B77").Value = "Vega Monumental Concepción"
$ws.Range("C77").Value = "Bíobío"
$ws.Range("D77").Value = 44272
$ws.Range("E77").Value = 8
$ws.Range("F77").Value = 100112028
$ws.Range("G77").Value = "Sandia"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 300
$ws.Range("K77").Value = 2500
$ws.Range("L77").Value = 2500
$ws.Range("M77").Value = 2500
$ws.Range("N77").Value = "$/unidad"
$ws.Range("O77").Value = "Región de O'Higgins"
$ws.Range("P77").Value = 2500
$ws.Range("Q77").Value = 1
$ws.Range("R77").Value = "Hortaliza"

# Row 78
$ws.Range("A78").Value = 11
$ws.Range("B78").Value = "Vega Monumental Concepción"
$ws.Range("C78").Value = "Bíobío"
$ws.Range("D78").Value = 44272
$ws.Range("E78").Value = 8
$ws.Range("F78").Value = 100112028
$ws.Range("G78").Value = "Sandia"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Segunda"
$ws.Range("J78").Value = 300
$ws.Range("K78").Value = 2000
$ws.Range("L78").Value = 2000
$ws.Range("M78").Value = 2000
$ws.Range("N78").Value = "$/unidad"
$ws.Range("O78").Value = "Región de O'Higgins"
$ws.Range("P78").Value = 2000
$ws.Range("Q78").Value = 1
$ws.Range("R78").Value = "Hortaliza"

# Apply the same date style (s="2") used on other D column cells to the new D75:D78 cells
$ws.Range("D71:D78").NumberFormat = "YYYY-MM-DD HH:MM:SS"
